# Excel COM-interop edit script ("excel update work, pt i")
#
# Target sheet: "CPs" (sheet3.xml). Refreshes the last three existing
# support-request rows (86-88) with a new batch of names/IDs and appends
# 68 brand-new rows (89-156) in the same shape, all dated 45584
# (2024-10-??) using a new custom "dd.mm.YYYY" number format. The E
# column formula is also updated to prefix the VLOOKUP result with the
# looked-up ID from column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CPs")

$data = @(
    @(86, 'Garrrfield', 35751242),
    @(87, 'LittleKnights', 55272769),
    @(88, 'Xorron', 46884042),
    @(89, 'Kyredneck30', 45717123),
    @(90, 'SoulKnightSK', 68482687),
    @(91, 'Schneeman', 52625609),
    @(92, 'Letsii', 46742185),
    @(93, 'Daeloan', 78221959),
    @(94, 'finanzamtt', 60805758),
    @(95, 'Sevenupurs', 60450559),
    @(96, 'ShaRopi69', 55798182),
    @(97, 'MR.Grinch', 54016906),
    @(98, 'xljhx31', 51309966),
    @(99, 'RJShuda', 44955250),
    @(100, 'Evita.Q', 42177593),
    @(101, 'DangerMouse', 41222530),
    @(102, 'dragSD', 40263474),
    @(103, 'Goldvale', 36918040),
    @(104, 'TrebleMaker', 36297525),
    @(105, 'KingLiz', 35924718),
    @(106, 'mido009', 41196673),
    @(107, 'Maria4612', 38974543),
    @(108, 'RaquelsHero', 36445384),
    @(109, 'I am Groot', 35973077),
    @(110, 'VON', 35437997),
    @(111, 'Alecks_', 70136709),
    @(112, 'itaca90909090', 54132050),
    @(113, 'Nefi85', 44311802),
    @(114, 'Permobil', 42022280),
    @(115, 'Niimphy', 35999933),
    @(116, 'Daut', 50140283),
    @(117, 'ZeroX_47', 42208168),
    @(118, 'Schockaletta', 36600286),
    @(119, 'LESINVINCIBLES', 31135467),
    @(120, 'POLAR-BEAR', 61121304),
    @(121, 'ShyDust', 60245369),
    @(122, 'SKIF3006', 49544841),
    @(123, 'NOAH', 47409011),
    @(124, 'Kevo1707', 47285414),
    @(125, 'Grizzly', 45802878),
    @(126, 'GhostSlayer', 47375039),
    @(127, 'Busijay', 45909959),
    @(128, 'Stefan333', 45637522),
    @(129, 'CharLee', 44549137),
    @(130, 'Ruby Sunday', 44237063),
    @(131, 'RockerFoo', 43864931),
    @(132, 'Hallen98', 43038395),
    @(133, 'The RoastPotato', 41356377),
    @(134, 'runyaover', 39811802),
    @(135, 'Luna_Lulu', 39057432),
    @(136, 'Lililulu', 384728894),
    @(137, 'OldManLogan', 37911832),
    @(138, 'VonTempsky', 37176660),
    @(139, 'Odin1206', 36118388),
    @(140, 'MizzBond', 35423556),
    @(141, 'Brilith', 38241597),
    @(142, 'Elo785', 37512673),
    @(143, 'Pilot_', 36293400),
    @(144, 'Lucky6612', 35810392),
    @(145, 'Omilixo', 34953211),
    @(146, 'chenjun', 33548915),
    @(147, 'Fab1250', 34546068),
    @(148, 'Valessa', 56755772),
    @(149, 'BasilFawlty', 75666591),
    @(150, 'LewkSkywatcher', 53752299),
    @(151, 'ocdMonkey', 54250384),
    @(152, 'DropKick4', 48578895),
    @(153, 'CMLTO', 66778954),
    @(154, 'WhoDey812', 59781609),
    @(155, 'Horizen', 60231598),
    @(156, 'Amcoone', 417634243)
)

$dateSerial = 45584
$dateFormat = "dd.mm.YYYY"

foreach ($entry in $data) {
    $r    = $entry[0]
    $name = $entry[1]
    $cid  = $entry[2]

    # B: reporter name, C: reporter numeric ID
    $ws.Cells.Item($r, 2).Value = $name
    $ws.Cells.Item($r, 3).Value = $cid

    # D: support-request date, formatted dd.mm.YYYY (new custom style)
    $ws.Cells.Item($r, 4).Value = $dateSerial
    $ws.Cells.Item($r, 4).NumberFormat = $dateFormat

    # A: ID lookup by name (unchanged formula shape)
    $ws.Cells.Item($r, 1).Formula = "=VLOOKUP(B$r,IDs!B:C,2,FALSE)"

    # E: now prefixes the looked-up support-list value with the ID from A
    $ws.Cells.Item($r, 5).Formula = "=A$r & ""|"" & VLOOKUP(D$r,SupportLists!D:E,2,FALSE)"
}

# The source edit also flips the workbook to recalc everything on next
# open (calcPr fullCalcOnLoad="1") since a batch of formulas was just
# dropped in/refreshed without a manual recalc pass.
$excel.Calculation = -4105
$excel.CalculateFullRebuild()
